# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" and "全部类型" sheets, which contain the same underlying rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    6  = 1351
    7  = 1584
    18 = 1773
    22 = 702
    25 = 4295
    31 = 650
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
